# Make Excel calculation/automation single-threaded (addresses random COM
# exceptions seen when running multi-threaded).
$excel.MultiThreadedCalculation = $false

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next trade record as row 6, mirroring the layout of the
# existing rows (Principle, Start Principle, BuyPrice, SellPrice,
# IsShortSell, Price Change %, Date, Profitable).
$ws.Range("A6").Value = 9849.31
$ws.Range("B6").Value = 9949.7999999999993
$ws.Range("C6").Value = 286
$ws.Range("D6").Value = 283.11
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = -1.01

# Column G carries a date-time style (style index 1). Copy the format from
# the cell above before setting its value so the new cell reuses the same
# style instead of Excel allocating a brand new one.
$ws.Range("G5").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("G6").Value = 42612.675057870372

$ws.Range("H6").Value = $false
